$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are always plain text; column E percentage strings never
# auto-parse as numbers. Column D sometimes holds plain numeric-looking text
# (e.g. "0.641", "39.20") which Excel would otherwise silently coerce into a
# real number (losing the original formatting / trailing zeros). For those
# cells we force the Text number format first so the literal string sticks.

$ws.Range("D2").Value = '42.946.82'
$ws.Range("E2").Value = '  +4.20%  '
$ws.Range("D3").Value = '2.284.50'
$ws.Range("E3").Value = '  +4.78%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.72'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.641'
$ws.Range("E6").Value = '  +4.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.72'
$ws.Range("E7").Value = '  +8.99%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.652'
$ws.Range("E9").Value = '  +12.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.20'
$ws.Range("E10").Value = '  +6.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0980'
$ws.Range("E11").Value = '  +4.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.82'
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.44'
$ws.Range("E13").Value = '  +8.17%  '
$ws.Range("E14").Value = '  +2.08%  '
$ws.Range("D15").Value = '2.626.51'
$ws.Range("E15").Value = '  +4.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.10'
$ws.Range("E16").Value = '  +5.44%  '
$ws.Range("E17").Value = '  +5.59%  '
$ws.Range("D18").Value = '2.289.85'
$ws.Range("E18").Value = '  +6.75%  '
$ws.Range("D19").Value = '42.880.28'
$ws.Range("E19").Value = '  +4.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000101'
$ws.Range("E20").Value = '  +6.83%  '
$ws.Range("E21").Value = '  +5.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.63'
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.80'
$ws.Range("E23").Value = '  +3.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.18'
$ws.Range("E24").Value = '  +8.07%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.83'
$ws.Range("E25").Value = '  +4.04%  '
$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.89'
$ws.Range("E26").Value = '  +2.09%  '
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.45'
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("E30").Value = '  +0.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.12'
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.16'
$ws.Range("E32").Value = '  +4.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.34'
$ws.Range("E33").Value = '  +11.75%  '
$ws.Range("E34").Value = '  +6.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0814'
$ws.Range("E35").Value = '  +8.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.45'
$ws.Range("E36").Value = '  +27.74%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.83'
$ws.Range("E37").Value = '  +22.25%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.127'
$ws.Range("E38").Value = '  +4.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.80'
$ws.Range("E39").Value = '  +6.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0311'
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.36'
$ws.Range("E41").Value = '  +17.06%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.33'
$ws.Range("E42").Value = '  +5.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.08'
$ws.Range("E43").Value = '  +10.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.213'
$ws.Range("E44").Value = '  +12.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.21'
$ws.Range("E45").Value = '  +8.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.96'
$ws.Range("E46").Value = '  -11.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '61.92'
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("E48").Value = '  +4.76%  '
$ws.Range("E49").Value = '  +4.84%  '
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  +5.12%  '

Write-Output "Applied 99 cell updates"
